$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts in row 2
$ws.Range("A2").Value = 33
$ws.Range("B2").Value = 29
$ws.Range("C2").Value = 3

# Update probabilities in row 5 (B2/A2, C2/A2, D2/A2)
$ws.Range("B5").Value = 0.8787878787878788
$ws.Range("C5").Value = 0.09090909090909091
$ws.Range("D5").Value = 0.0303030303030303
